# Updates cryptos list (Price/Volume(1h) columns) to the latest scraped
# values, including two rows whose coins were reordered (Monero/Cosmos at
# rows 26-27 and ARBITRUM/HuobiToken at rows 47-48).
#
# Price values (column D) that look like plain numbers are written with a
# leading apostrophe so Excel keeps them as literal text (e.g. "227.90")
# instead of coercing them into a numeric value (227.9) and losing the
# original formatting/trailing zero.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.652.78'
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').Value = '2.163.02'
$ws.Range('E3').Value = '  +2.76%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '''227.90'
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('D6').Value = '''0.631'
$ws.Range('E6').Value = '  +2.36%  '
$ws.Range('D7').Value = '''63.61'
$ws.Range('E7').Value = '  +2.12%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').Value = '''0.393'
$ws.Range('E9').Value = '  +0.77%  '
$ws.Range('E10').Value = '  +0.67%  '
$ws.Range('E11').Value = '  +0.04%  '
$ws.Range('D12').Value = '''16.04'
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('D13').Value = '2.484.46'
$ws.Range('E13').Value = '  +2.79%  '
$ws.Range('D14').Value = '''22.01'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('D16').Value = '''5.49'
$ws.Range('E16').Value = '  -0.76%  '
$ws.Range('D17').Value = '2.157.53'
$ws.Range('E17').Value = '  +4.01%  '
$ws.Range('D18').Value = '39.591.26'
$ws.Range('E18').Value = '  +1.98%  '
$ws.Range('D19').Value = '''71.86'
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('D20').Value = '''6.12'
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').Value = '0.0₃0846'
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('D22').Value = '''228.06'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').Value = '''2.41'
$ws.Range('E24').Value = '  +3.56%  '
$ws.Range('D25').Value = '''2.37'
$ws.Range('E25').Value = '  +0.25%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '''9.66'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '''172.43'
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('E28').Value = '  +0.88%  '
$ws.Range('D29').Value = '''19.75'
$ws.Range('E29').Value = '  +2.08%  '
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('E31').Value = '  +4.44%  '
$ws.Range('E32').Value = '  +1.59%  '
$ws.Range('D33').Value = '''4.60'
$ws.Range('E33').Value = '  +0.54%  '
$ws.Range('E34').Value = '  -1.47%  '
$ws.Range('D35').Value = '''6.96'
$ws.Range('E35').Value = '  -3.01%  '
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('E37').Value = '  +0.63%  '
$ws.Range('E38').Value = '  +3.01%  '
$ws.Range('D39').Value = '''1.00'
$ws.Range('E39').Value = '  +0.45%  '
$ws.Range('D40').Value = '''4.87'
$ws.Range('E40').Value = '  +16.88%  '
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('E43').Value = '  -2.58%  '
$ws.Range('D44').Value = '''1.24'
$ws.Range('E44').Value = '  +2.23%  '
$ws.Range('D45').Value = '1.513.39'
$ws.Range('E45').Value = '  -0.81%  '
$ws.Range('D46').Value = '''0.0924'
$ws.Range('E46').Value = '  +0.69%  '
$ws.Range('B47').Value = 'HuobiToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D47').Value = '''2.80'
$ws.Range('E47').Value = '  -0.23%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').Value = '''1.10'
$ws.Range('E48').Value = '  +1.50%  '
$ws.Range('D49').Value = '''7.78'
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('E50').Value = '  +1.26%  '
$ws.Range('D51').Value = '2.368.87'
$ws.Range('E51').Value = '  +2.80%  '
